$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 0.4165321785010008
$ws.Range("D2").Value = 0.1770283408321944
$ws.Range("E2").Value = 0.1759721268985714
$ws.Range("F2").Value = 3.387674600170385
$ws.Range("G2").Value = 3.717185640561695
$ws.Range("H2").Value = 2.478075300866522
$ws.Range("L2").Value = 0.1535809029467714
$ws.Range("N2").Value = 2.873581687696401
$ws.Range("C3").Value = 0.4086178312782067
$ws.Range("D3").Value = 0.1696488696270393
$ws.Range("E3").Value = 0.1722178120729154
$ws.Range("F3").Value = 3.208307731575303
$ws.Range("G3").Value = 3.493995001652593
$ws.Range("H3").Value = 2.383840216937983
$ws.Range("L3").Value = 0.1499384042476493
$ws.Range("N3").Value = 2.562605684679681
$ws.Range("C4").Value = 0.4040558284518738
$ws.Range("D4").Value = 0.1652659502853737
$ws.Range("E4").Value = 0.1700421457223484
$ws.Range("F4").Value = 3.10095405640476
$ws.Range("G4").Value = 3.360036834047719
$ws.Range("H4").Value = 2.32794307985273
$ws.Range("L4").Value = 0.1478181069540341
$ws.Range("N4").Value = 2.371325805375818
$ws.Range("C5").Value = 0.4022706003519829
$ws.Range("D5").Value = 0.1635162936334638
$ws.Range("E5").Value = 0.1691876484232147
$ws.Range("F5").Value = 3.057888776803537
$ws.Range("G5").Value = 3.306202523924014
$ws.Range("H5").Value = 2.305649011294406
$ws.Range("L5").Value = 0.1469828456313849
$ws.Range("N5").Value = 2.293303068605894
$ws.Range("C6").Value = 0.4019785956468525
$ws.Range("D6").Value = 0.1632279403856387
$ws.Range("E6").Value = 0.1690476856441485
$ws.Range("F6").Value = 3.050778551949776
$ws.Range("G6").Value = 3.297308380218112
$ws.Range("H6").Value = 2.301976091865981
$ws.Range("L6").Value = 0.1468458763935558
$ws.Range("N6").Value = 2.2803432614038
$ws.Range("C7").Value = 0.4040314545854926
$ws.Range("D7").Value = 0.165242207419297
$ws.Range("E7").Value = 0.1700304922803255
$ws.Range("F7").Value = 3.100370522109301
$ws.Range("G7").Value = 3.359307776180344
$ws.Range("H7").Value = 2.32764046474864
$ws.Range("L7").Value = 0.1478067263900584
$ws.Range("N7").Value = 2.370273851392596
$ws.Range("C8").Value = 0.4137409840663224
$ws.Range("D8").Value = 0.1744526551100023
$ws.Range("E8").Value = 0.1746504803007056
$ws.Range("F8").Value = 3.325242206422672
$ws.Range("G8").Value = 3.639577500156065
$ws.Range("H8").Value = 2.4451698046758
$ws.Range("L8").Value = 0.1523005766433769
$ws.Range("N8").Value = 2.766433886209882
$ws.Range("C9").Value = 0.4351876940654904
$ws.Range("D9").Value = 0.1937285989377386
$ws.Range("E9").Value = 0.1847596816637775
$ws.Range("F9").Value = 3.789042000569282
$ws.Range("G9").Value = 4.21460069015319
$ws.Range("H9").Value = 2.69166815033725
$ws.Range("L9").Value = 0.1620563131879749
$ws.Range("N9").Value = 3.540180268007646
$ws.Range("C10").Value = 0.452479671672819
$ws.Range("D10").Value = 0.208688054512379
$ws.Range("E10").Value = 0.1928588232532249
$ws.Range("F10").Value = 4.144858485178389
$ws.Range("G10").Value = 4.653973240002585
$ws.Range("H10").Value = 2.883189741785998
$ws.Range("L10").Value = 0.1698300199493019
$ws.Range("N10").Value = 4.10623028343673
$ws.Range("C11").Value = 0.4606950994645729
$ws.Range("D11").Value = 0.2156796332962756
$ws.Range("E11").Value = 0.1966965571950539
$ws.Range("F11").Value = 4.310259620081865
$ws.Range("G11").Value = 4.857845151148183
$ws.Range("H11").Value = 2.97272826291487
$ws.Range("L11").Value = 0.1735051820266733
$ws.Range("N11").Value = 4.363110593465422
$ws.Range("C12").Value = 0.4638575714749607
$ws.Range("D12").Value = 0.2183550492997881
$ws.Range("E12").Value = 0.1981724717683875
$ws.Range("F12").Value = 4.373422961433562
$ws.Range("G12").Value = 4.935647507225099
$ws.Range("H12").Value = 3.006993552294546
$ws.Range("L12").Value = 0.1749174252542218
$ws.Range("N12").Value = 4.460285735713398
$ws.Range("C13").Value = 0.4631741670923759
$ws.Range("D13").Value = 0.2177775946004488
$ws.Range("E13").Value = 0.1978535902517393
$ws.Range("F13").Value = 4.35979571788647
$ws.Range("G13").Value = 4.918864256874883
$ws.Range("H13").Value = 2.999597750250189
$ws.Range("L13").Value = 0.174612350616286
$ws.Range("N13").Value = 4.439361943450422
$ws.Range("C14").Value = 0.460954239117882
$ws.Range("D14").Value = 0.2158991769628358
$ws.Range("E14").Value = 0.1968175242162289
$ws.Range("F14").Value = 4.315445372688544
$ws.Range("G14").Value = 4.864233814429042
$ws.Range("H14").Value = 2.975540024997656
$ws.Range("L14").Value = 0.1736209531488981
$ws.Range("N14").Value = 4.371107314139294
$ws.Range("C15").Value = 0.4596012096658626
$ws.Range("D15").Value = 0.2147522513856472
$ws.Range("E15").Value = 0.1961858711810578
$ws.Range("F15").Value = 4.288349091472526
$ws.Range("G15").Value = 4.830850061484398
$ws.Range("H15").Value = 2.960851082899751
$ws.Range("L15").Value = 0.1730163856637148
$ws.Range("N15").Value = 4.329286057409945
$ws.Range("C16").Value = 0.4519499071874975
$ws.Range("D16").Value = 0.2082349787651481
$ws.Range("E16").Value = 0.1926111549514715
$ws.Range("F16").Value = 4.134122087581517
$ws.Range("G16").Value = 4.64073235703637
$ws.Range("H16").Value = 2.87738777250263
$ws.Range("L16").Value = 0.1695926820739118
$ws.Range("N16").Value = 4.089429168003562
$ws.Range("C17").Value = 0.4473463998459408
$ws.Range("D17").Value = 0.2042853293912401
$ws.Range("E17").Value = 0.1904578833456796
$ws.Range("F17").Value = 4.040429598389494
$ws.Range("G17").Value = 4.525143391300105
$ws.Range("H17").Value = 2.826812389449003
$ws.Range("L17").Value = 0.1675283167678998
$ws.Range("N17").Value = 3.94211849063862
$ws.Range("C18").Value = 0.4447313769314292
$ws.Range("D18").Value = 0.2020310800897107
$ws.Range("E18").Value = 0.1892337800419384
$ws.Range("F18").Value = 3.986871930230961
$ws.Range("G18").Value = 4.459034521134527
$ws.Range("H18").Value = 2.797949289413225
$ws.Range("L18").Value = 0.1663539836133054
$ws.Range("N18").Value = 3.857331695637129
$ws.Range("C19").Value = 0.4438515699405912
$ws.Range("D19").Value = 0.2012708017276452
$ws.Range("E19").Value = 0.1888217755454988
$ws.Range("F19").Value = 3.968794587278524
$ws.Range("G19").Value = 4.43671481326902
$ws.Range("H19").Value = 2.788215292087557
$ws.Range("L19").Value = 0.1659585963334393
$ws.Range("N19").Value = 3.828614786364199
$ws.Range("C20").Value = 0.447833048226471
$ws.Range("D20").Value = 0.2047039592653448
$ws.Range("E20").Value = 0.1906856080695505
$ws.Range("F20").Value = 4.050368838828405
$ws.Range("G20").Value = 4.537409057556317
$ws.Range("H20").Value = 2.832172696867701
$ws.Range("L20").Value = 0.1677467188202399
$ws.Range("N20").Value = 3.957806003277312
$ws.Range("C21").Value = 0.4616048787614488
$ws.Range("D21").Value = 0.2164501492096917
$ws.Range("E21").Value = 0.1971212224547401
$ws.Range("F21").Value = 4.32845760456712
$ws.Range("G21").Value = 4.880263595686813
$ws.Range("H21").Value = 2.98259652285725
$ws.Range("L21").Value = 0.1739115885413014
$ws.Range("N21").Value = 4.391158149571254
$ws.Range("C22").Value = 0.4709061736978981
$ws.Range("D22").Value = 0.2242898467258669
$ws.Range("E22").Value = 0.2014595680060083
$ws.Range("F22").Value = 4.513301266596727
$ws.Range("G22").Value = 5.107851384151218
$ws.Range("H22").Value = 3.083005467025714
$ws.Range("L22").Value = 0.1780606847380568
$ws.Range("N22").Value = 4.673791817957863
$ws.Range("C23").Value = 0.4659139666193539
$ws.Range("D23").Value = 0.2200903934251812
$ws.Range("E23").Value = 0.1991318052414854
$ws.Range("F23").Value = 4.414356361456555
$ws.Range("G23").Value = 4.986053475071742
$ws.Range("H23").Value = 3.02921931974447
$ws.Range("L23").Value = 0.1758350602497245
$ws.Range("N23").Value = 4.523002190001307
$ws.Range("C24").Value = 0.447612936177336
$ws.Range("D24").Value = 0.204514645519339
$ws.Range("E24").Value = 0.1905826106539976
$ws.Range("F24").Value = 4.04587435232574
$ws.Range("G24").Value = 4.531862677884931
$ws.Range("H24").Value = 2.829748638198907
$ws.Range("L24").Value = 0.1676479403671038
$ws.Range("N24").Value = 3.950713976768498
$ws.Range("C25").Value = 0.4291212773186999
$ws.Range("D25").Value = 0.1883781532025068
$ws.Range("E25").Value = 0.1819092697119515
$ws.Range("F25").Value = 3.661011727456582
$ws.Range("G25").Value = 4.056177817182629
$ws.Range("H25").Value = 2.623202387292906
$ws.Range("L25").Value = 0.1593129888777014
$ws.Range("N25").Value = 3.331249627311138
